# Updated symbol list (coinranking.com scrape refresh) — GitHub Actions run.
# Each row's Price(D)/rank-prefixed Volume(E)/Hour(G) gets refreshed; the
# coin roster in rows 10-24 shifted down one slot to make room for a new
# "One" (ONE) entry pulled in at the top of that block.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell => new value. "Numeric" cells hold digit-only text (prices, hours)
# that Excel would otherwise auto-coerce to a float and mangle trailing
# zeros (e.g. "0.01120" -> 0.0112), so we force the Text number format
# before writing those.
$updates = @(
    @{ Cell = "D2"; Value = "246.11"; Numeric = $true },
    @{ Cell = "G2"; Value = "19"; Numeric = $true },
    @{ Cell = "D3"; Value = "22.24"; Numeric = $true },
    @{ Cell = "G3"; Value = "19"; Numeric = $true },
    @{ Cell = "G4"; Value = "19"; Numeric = $true },
    @{ Cell = "D5"; Value = "0.05859"; Numeric = $true },
    @{ Cell = "G5"; Value = "19"; Numeric = $true },
    @{ Cell = "G6"; Value = "19"; Numeric = $true },
    @{ Cell = "D7"; Value = "6.380"; Numeric = $true },
    @{ Cell = "G7"; Value = "19"; Numeric = $true },
    @{ Cell = "D8"; Value = "0.8129"; Numeric = $true },
    @{ Cell = "G8"; Value = "19"; Numeric = $true },
    @{ Cell = "D9"; Value = "0.9928"; Numeric = $true },
    @{ Cell = "G9"; Value = "19"; Numeric = $true },
    @{ Cell = "B10"; Value = "One"; Numeric = $false },
    @{ Cell = "C10"; Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"; Numeric = $false },
    @{ Cell = "D10"; Value = "0.01120"; Numeric = $true },
    @{ Cell = "E10"; Value = "9OneONEBestin24h"; Numeric = $false },
    @{ Cell = "G10"; Value = "19"; Numeric = $true },
    @{ Cell = "B11"; Value = "WazirX"; Numeric = $false },
    @{ Cell = "C11"; Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"; Numeric = $false },
    @{ Cell = "D11"; Value = "0.1422"; Numeric = $true },
    @{ Cell = "E11"; Value = "10WazirXWRX"; Numeric = $false },
    @{ Cell = "G11"; Value = "19"; Numeric = $true },
    @{ Cell = "B12"; Value = "LiechtensteinCryptoassetsExchange"; Numeric = $false },
    @{ Cell = "C12"; Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"; Numeric = $false },
    @{ Cell = "D12"; Value = "0.03987"; Numeric = $true },
    @{ Cell = "E12"; Value = "11LiechtensteinCryptoassetsExchangeLCX"; Numeric = $false },
    @{ Cell = "G12"; Value = "19"; Numeric = $true },
    @{ Cell = "B13"; Value = "MandalaExchangeToken"; Numeric = $false },
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"; Numeric = $false },
    @{ Cell = "D13"; Value = "0.07344"; Numeric = $true },
    @{ Cell = "E13"; Value = "12MandalaExchangeTokenMDX"; Numeric = $false },
    @{ Cell = "G13"; Value = "19"; Numeric = $true },
    @{ Cell = "B14"; Value = "BitrueCoin"; Numeric = $false },
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"; Numeric = $false },
    @{ Cell = "D14"; Value = "0.03007"; Numeric = $true },
    @{ Cell = "E14"; Value = "13BitrueCoinBTR"; Numeric = $false },
    @{ Cell = "G14"; Value = "19"; Numeric = $true },
    @{ Cell = "B15"; Value = "MCDex"; Numeric = $false },
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"; Numeric = $false },
    @{ Cell = "D15"; Value = "4.178"; Numeric = $true },
    @{ Cell = "E15"; Value = "14MCDexMCB"; Numeric = $false },
    @{ Cell = "G15"; Value = "19"; Numeric = $true },
    @{ Cell = "B16"; Value = "BitMartToken"; Numeric = $false },
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"; Numeric = $false },
    @{ Cell = "D16"; Value = "0.09392"; Numeric = $true },
    @{ Cell = "E16"; Value = "15BitMartTokenBMX"; Numeric = $false },
    @{ Cell = "G16"; Value = "19"; Numeric = $true },
    @{ Cell = "B17"; Value = "BitForexToken"; Numeric = $false },
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"; Numeric = $false },
    @{ Cell = "D17"; Value = "0.001591"; Numeric = $true },
    @{ Cell = "E17"; Value = "16BitForexTokenBF"; Numeric = $false },
    @{ Cell = "G17"; Value = "19"; Numeric = $true },
    @{ Cell = "B18"; Value = "CoinExToken"; Numeric = $false },
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"; Numeric = $false },
    @{ Cell = "D18"; Value = "0.04809"; Numeric = $true },
    @{ Cell = "E18"; Value = "17CoinExTokenCET"; Numeric = $false },
    @{ Cell = "G18"; Value = "19"; Numeric = $true },
    @{ Cell = "B19"; Value = "TigerCash"; Numeric = $false },
    @{ Cell = "C19"; Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"; Numeric = $false },
    @{ Cell = "D19"; Value = "0.006102"; Numeric = $true },
    @{ Cell = "E19"; Value = "18TigerCashTCH"; Numeric = $false },
    @{ Cell = "G19"; Value = "19"; Numeric = $true },
    @{ Cell = "B20"; Value = "HotbitToken"; Numeric = $false },
    @{ Cell = "C20"; Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"; Numeric = $false },
    @{ Cell = "D20"; Value = "0.004080"; Numeric = $true },
    @{ Cell = "E20"; Value = "19HotbitTokenHTB"; Numeric = $false },
    @{ Cell = "G20"; Value = "19"; Numeric = $true },
    @{ Cell = "B21"; Value = "BitKan"; Numeric = $false },
    @{ Cell = "C21"; Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"; Numeric = $false },
    @{ Cell = "D21"; Value = "0.0009839"; Numeric = $true },
    @{ Cell = "E21"; Value = "20BitKanKAN"; Numeric = $false },
    @{ Cell = "G21"; Value = "19"; Numeric = $true },
    @{ Cell = "B22"; Value = "NitroEx"; Numeric = $false },
    @{ Cell = "C22"; Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"; Numeric = $false },
    @{ Cell = "D22"; Value = "0.0001410"; Numeric = $true },
    @{ Cell = "E22"; Value = "21NitroExNTX"; Numeric = $false },
    @{ Cell = "G22"; Value = "19"; Numeric = $true },
    @{ Cell = "B23"; Value = "LEO"; Numeric = $false },
    @{ Cell = "C23"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; Numeric = $false },
    @{ Cell = "D23"; Value = "3.690"; Numeric = $true },
    @{ Cell = "E23"; Value = "22LEOLEO"; Numeric = $false },
    @{ Cell = "G23"; Value = "19"; Numeric = $true },
    @{ Cell = "B24"; Value = "BTSEToken"; Numeric = $false },
    @{ Cell = "C24"; Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"; Numeric = $false },
    @{ Cell = "D24"; Value = "2.233"; Numeric = $true },
    @{ Cell = "E24"; Value = "23BTSETokenBTSE"; Numeric = $false },
    @{ Cell = "G24"; Value = "19"; Numeric = $true },
    @{ Cell = "D25"; Value = "0.3246"; Numeric = $true },
    @{ Cell = "G25"; Value = "19"; Numeric = $true },
    @{ Cell = "D26"; Value = "0.1295"; Numeric = $true },
    @{ Cell = "G26"; Value = "19"; Numeric = $true },
    @{ Cell = "G27"; Value = "19"; Numeric = $true },
    @{ Cell = "G28"; Value = "19"; Numeric = $true },
    @{ Cell = "G29"; Value = "19"; Numeric = $true },
    @{ Cell = "G30"; Value = "19"; Numeric = $true },
    @{ Cell = "G31"; Value = "19"; Numeric = $true },
    @{ Cell = "G32"; Value = "19"; Numeric = $true },
    @{ Cell = "G33"; Value = "19"; Numeric = $true },
    @{ Cell = "G34"; Value = "19"; Numeric = $true },
    @{ Cell = "G35"; Value = "19"; Numeric = $true },
    @{ Cell = "G36"; Value = "19"; Numeric = $true },
    @{ Cell = "G37"; Value = "19"; Numeric = $true },
    @{ Cell = "G38"; Value = "19"; Numeric = $true },
    @{ Cell = "G39"; Value = "19"; Numeric = $true },
    @{ Cell = "D40"; Value = "0.03858"; Numeric = $true },
    @{ Cell = "G40"; Value = "19"; Numeric = $true },
    @{ Cell = "D41"; Value = "0.006401"; Numeric = $true },
    @{ Cell = "G41"; Value = "19"; Numeric = $true },
    @{ Cell = "D42"; Value = "0.1073"; Numeric = $true },
    @{ Cell = "G42"; Value = "19"; Numeric = $true },
    @{ Cell = "D43"; Value = "0.002601"; Numeric = $true },
    @{ Cell = "G43"; Value = "19"; Numeric = $true },
    @{ Cell = "D44"; Value = "0.005207"; Numeric = $true },
    @{ Cell = "G44"; Value = "19"; Numeric = $true },
    @{ Cell = "D45"; Value = "0.00005650"; Numeric = $true },
    @{ Cell = "G45"; Value = "19"; Numeric = $true },
    @{ Cell = "G46"; Value = "19"; Numeric = $true },
    @{ Cell = "D47"; Value = "0.7222"; Numeric = $true },
    @{ Cell = "G47"; Value = "19"; Numeric = $true },
    @{ Cell = "D48"; Value = "0.08619"; Numeric = $true },
    @{ Cell = "G48"; Value = "19"; Numeric = $true },
    @{ Cell = "D49"; Value = "0.00002101"; Numeric = $true },
    @{ Cell = "G49"; Value = "19"; Numeric = $true },
    @{ Cell = "G50"; Value = "19"; Numeric = $true },
    @{ Cell = "G51"; Value = "19"; Numeric = $true }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Numeric) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}

Write-Host "Applied $($updates.Count) cell updates"
